$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")
$ws.Activate()

# Rows 3 and 4 (the "Immersive modeling" and "Capped depletions" studies)
# need to swap places in the table - read current values first.
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2

$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$d4 = $ws.Range("D4").Value2
$e4 = $ws.Range("E4").Value2

$ws.Range("A3").Value = $a4
$ws.Range("B3").Value = $b4
$ws.Range("D3").Value = $d4
$ws.Range("E3").Value = $e4

$ws.Range("A4").Value = $a3
$ws.Range("B4").Value = $b3
$ws.Range("D4").Value = $d3
$ws.Range("E4").Value = $e3

# Update the selection left behind on the sheet to cover the entire row 3
$ws.Range("A3:XFD3").Select() | Out-Null
